$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the cell formatting (style) used by columns A:B (rows 2-10) onto
# columns C:D for the same rows, so the new data mirrors the existing look.
$ws.Range("A2:B10").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill columns C and D (rows 2-10) with the same values already present in
# columns A and B respectively.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 2).Value2
}

# Update the active selection to a single cell.
$ws.Range("C8").Select() | Out-Null
